# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") held the old "Strike#" derived values. The data
# pipeline was regenerated so that column G now holds the recomputed
# "K" (strikeout-based) values. The new values for each row are written
# below, keyed by worksheet row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 2
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 3
    11 = 2
    12 = 1
    13 = 2
    14 = 5
    15 = 3
    16 = 5
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    31 = 0
    32 = 2
    33 = 0
    34 = 2
    35 = 0
    36 = 3
    37 = 2
    38 = 2
    39 = 2
    40 = 1
    41 = 2
    42 = 1
    43 = 0
    44 = 0
    45 = 1
    46 = 2
    47 = 2
    48 = 2
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 2
    54 = 2
    55 = 1
    56 = 1
    57 = 1
    58 = 1
    59 = 3
    60 = 3
    61 = 1
    62 = 1
    63 = 2
    64 = 2
    65 = 0
    66 = 2
    67 = 3
    68 = 0
    69 = 1
    70 = 3
    71 = 0
    72 = 0
    73 = 2
    74 = 2
    75 = 1
    76 = 0
    77 = 1
    78 = 1
    79 = 1
    80 = 2
    81 = 1
    83 = 0
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
